$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 4 changed model-name labels (column A) for each block of 10 rows
$ws.Range("A2:A11").Value = "CNN_Attention_32_filters_5_kernels_predictions"
$ws.Range("A12:A21").Value = "CNN_32_filters_3_kernels_predictions"
$ws.Range("A22:A31").Value = "CNN_64_filters_7_kernels_predictions"
$ws.Range("A32:A41").Value = "CNN_Attention_256_filters_5_kernels_predictions"

# Update numeric metric cells E:T for rows 2-41 (cell by cell, since bulk array assignment is not supported)
# Row 2
$ws.Cells.Item(2, 5).Value = 0.7818181818181819
$ws.Cells.Item(2, 6).Value = 0.7413793103448276
$ws.Cells.Item(2, 7).Value = 0.7610619469026548
$ws.Cells.Item(2, 8).Value = 0.7818181818181819
$ws.Cells.Item(2, 9).Value = 0.7413793103448276
$ws.Cells.Item(2, 10).Value = 0.7610619469026548
$ws.Cells.Item(2, 11).Value = 0.1563342318059299
$ws.Cells.Item(2, 12).Value = 0.7894736842105263
$ws.Cells.Item(2, 13).Value = 0.7758620689655172
$ws.Cells.Item(2, 14).Value = 0.7826086956521741
$ws.Cells.Item(2, 15).Value = 0.007655502392344471
$ws.Cells.Item(2, 16).Value = 0.03448275862068961
$ws.Cells.Item(2, 17).Value = 0.02154674874951923
$ws.Cells.Item(2, 18).Value = 0.009791921664626649
$ws.Cells.Item(2, 19).Value = 0.04651162790697669
$ws.Cells.Item(2, 20).Value = 0.02831142568250783
# Row 3
$ws.Cells.Item(3, 5).Value = 0.7818181818181819
$ws.Cells.Item(3, 6).Value = 0.7413793103448276
$ws.Cells.Item(3, 7).Value = 0.7610619469026548
$ws.Cells.Item(3, 8).Value = 0.7818181818181819
$ws.Cells.Item(3, 9).Value = 0.7413793103448276
$ws.Cells.Item(3, 10).Value = 0.7610619469026548
$ws.Cells.Item(3, 11).Value = 0.1563342318059299
$ws.Cells.Item(3, 12).Value = 0.7894736842105263
$ws.Cells.Item(3, 13).Value = 0.7758620689655172
$ws.Cells.Item(3, 14).Value = 0.7826086956521741
$ws.Cells.Item(3, 15).Value = 0.007655502392344471
$ws.Cells.Item(3, 16).Value = 0.03448275862068961
$ws.Cells.Item(3, 17).Value = 0.02154674874951923
$ws.Cells.Item(3, 18).Value = 0.009791921664626649
$ws.Cells.Item(3, 19).Value = 0.04651162790697669
$ws.Cells.Item(3, 20).Value = 0.02831142568250783
# Row 4
$ws.Cells.Item(4, 5).Value = 0.7818181818181819
$ws.Cells.Item(4, 6).Value = 0.7413793103448276
$ws.Cells.Item(4, 7).Value = 0.7610619469026548
$ws.Cells.Item(4, 8).Value = 0.7818181818181819
$ws.Cells.Item(4, 9).Value = 0.7413793103448276
$ws.Cells.Item(4, 10).Value = 0.7610619469026548
$ws.Cells.Item(4, 11).Value = 0.1563342318059299
$ws.Cells.Item(4, 12).Value = 0.7818181818181819
$ws.Cells.Item(4, 13).Value = 0.7413793103448276
$ws.Cells.Item(4, 14).Value = 0.7610619469026548
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
# Row 5
$ws.Cells.Item(5, 5).Value = 0.7818181818181819
$ws.Cells.Item(5, 6).Value = 0.7413793103448276
$ws.Cells.Item(5, 7).Value = 0.7610619469026548
$ws.Cells.Item(5, 8).Value = 0.7818181818181819
$ws.Cells.Item(5, 9).Value = 0.7413793103448276
$ws.Cells.Item(5, 10).Value = 0.7610619469026548
$ws.Cells.Item(5, 11).Value = 0.1563342318059299
$ws.Cells.Item(5, 12).Value = 0.7894736842105263
$ws.Cells.Item(5, 13).Value = 0.7758620689655172
$ws.Cells.Item(5, 14).Value = 0.7826086956521741
$ws.Cells.Item(5, 15).Value = 0.007655502392344471
$ws.Cells.Item(5, 16).Value = 0.03448275862068961
$ws.Cells.Item(5, 17).Value = 0.02154674874951923
$ws.Cells.Item(5, 18).Value = 0.009791921664626649
$ws.Cells.Item(5, 19).Value = 0.04651162790697669
$ws.Cells.Item(5, 20).Value = 0.02831142568250783
# Row 6
$ws.Cells.Item(6, 5).Value = 0.7818181818181819
$ws.Cells.Item(6, 6).Value = 0.7413793103448276
$ws.Cells.Item(6, 7).Value = 0.7610619469026548
$ws.Cells.Item(6, 8).Value = 0.7818181818181819
$ws.Cells.Item(6, 9).Value = 0.7413793103448276
$ws.Cells.Item(6, 10).Value = 0.7610619469026548
$ws.Cells.Item(6, 11).Value = 0.1563342318059299
$ws.Cells.Item(6, 12).Value = 0.7894736842105263
$ws.Cells.Item(6, 13).Value = 0.7758620689655172
$ws.Cells.Item(6, 14).Value = 0.7826086956521741
$ws.Cells.Item(6, 15).Value = 0.007655502392344471
$ws.Cells.Item(6, 16).Value = 0.03448275862068961
$ws.Cells.Item(6, 17).Value = 0.02154674874951923
$ws.Cells.Item(6, 18).Value = 0.009791921664626649
$ws.Cells.Item(6, 19).Value = 0.04651162790697669
$ws.Cells.Item(6, 20).Value = 0.02831142568250783
# Row 7
$ws.Cells.Item(7, 5).Value = 0.7818181818181819
$ws.Cells.Item(7, 6).Value = 0.7413793103448276
$ws.Cells.Item(7, 7).Value = 0.7610619469026548
$ws.Cells.Item(7, 8).Value = 0.7818181818181819
$ws.Cells.Item(7, 9).Value = 0.7413793103448276
$ws.Cells.Item(7, 10).Value = 0.7610619469026548
$ws.Cells.Item(7, 11).Value = 0.1563342318059299
$ws.Cells.Item(7, 12).Value = 0.7857142857142857
$ws.Cells.Item(7, 13).Value = 0.7586206896551724
$ws.Cells.Item(7, 14).Value = 0.7719298245614034
$ws.Cells.Item(7, 15).Value = 0.003896103896103842
$ws.Cells.Item(7, 16).Value = 0.01724137931034475
$ws.Cells.Item(7, 17).Value = 0.01086787765874853
$ws.Cells.Item(7, 18).Value = 0.004983388704318867
$ws.Cells.Item(7, 19).Value = 0.02325581395348827
$ws.Cells.Item(7, 20).Value = 0.01427988576091376
# Row 8
$ws.Cells.Item(8, 5).Value = 0.7818181818181819
$ws.Cells.Item(8, 6).Value = 0.7413793103448276
$ws.Cells.Item(8, 7).Value = 0.7610619469026548
$ws.Cells.Item(8, 8).Value = 0.7818181818181819
$ws.Cells.Item(8, 9).Value = 0.7413793103448276
$ws.Cells.Item(8, 10).Value = 0.7610619469026548
$ws.Cells.Item(8, 11).Value = 0.1563342318059299
$ws.Cells.Item(8, 12).Value = 0.7857142857142857
$ws.Cells.Item(8, 13).Value = 0.7586206896551724
$ws.Cells.Item(8, 14).Value = 0.7719298245614034
$ws.Cells.Item(8, 15).Value = 0.003896103896103842
$ws.Cells.Item(8, 16).Value = 0.01724137931034475
$ws.Cells.Item(8, 17).Value = 0.01086787765874853
$ws.Cells.Item(8, 18).Value = 0.004983388704318867
$ws.Cells.Item(8, 19).Value = 0.02325581395348827
$ws.Cells.Item(8, 20).Value = 0.01427988576091376
# Row 9
$ws.Cells.Item(9, 5).Value = 0.7818181818181819
$ws.Cells.Item(9, 6).Value = 0.7413793103448276
$ws.Cells.Item(9, 7).Value = 0.7610619469026548
$ws.Cells.Item(9, 8).Value = 0.7818181818181819
$ws.Cells.Item(9, 9).Value = 0.7413793103448276
$ws.Cells.Item(9, 10).Value = 0.7610619469026548
$ws.Cells.Item(9, 11).Value = 0.1563342318059299
$ws.Cells.Item(9, 12).Value = 0.7857142857142857
$ws.Cells.Item(9, 13).Value = 0.7586206896551724
$ws.Cells.Item(9, 14).Value = 0.7719298245614034
$ws.Cells.Item(9, 15).Value = 0.003896103896103842
$ws.Cells.Item(9, 16).Value = 0.01724137931034475
$ws.Cells.Item(9, 17).Value = 0.01086787765874853
$ws.Cells.Item(9, 18).Value = 0.004983388704318867
$ws.Cells.Item(9, 19).Value = 0.02325581395348827
$ws.Cells.Item(9, 20).Value = 0.01427988576091376
# Row 10
$ws.Cells.Item(10, 5).Value = 0.7818181818181819
$ws.Cells.Item(10, 6).Value = 0.7413793103448276
$ws.Cells.Item(10, 7).Value = 0.7610619469026548
$ws.Cells.Item(10, 8).Value = 0.7818181818181819
$ws.Cells.Item(10, 9).Value = 0.7413793103448276
$ws.Cells.Item(10, 10).Value = 0.7610619469026548
$ws.Cells.Item(10, 11).Value = 0.1563342318059299
$ws.Cells.Item(10, 12).Value = 0.7857142857142857
$ws.Cells.Item(10, 13).Value = 0.7586206896551724
$ws.Cells.Item(10, 14).Value = 0.7719298245614034
$ws.Cells.Item(10, 15).Value = 0.003896103896103842
$ws.Cells.Item(10, 16).Value = 0.01724137931034475
$ws.Cells.Item(10, 17).Value = 0.01086787765874853
$ws.Cells.Item(10, 18).Value = 0.004983388704318867
$ws.Cells.Item(10, 19).Value = 0.02325581395348827
$ws.Cells.Item(10, 20).Value = 0.01427988576091376
# Row 11
$ws.Cells.Item(11, 5).Value = 0.7818181818181819
$ws.Cells.Item(11, 6).Value = 0.7413793103448276
$ws.Cells.Item(11, 7).Value = 0.7610619469026548
$ws.Cells.Item(11, 8).Value = 0.7818181818181819
$ws.Cells.Item(11, 9).Value = 0.7413793103448276
$ws.Cells.Item(11, 10).Value = 0.7610619469026548
$ws.Cells.Item(11, 11).Value = 0.1563342318059299
$ws.Cells.Item(11, 12).Value = 0.7857142857142857
$ws.Cells.Item(11, 13).Value = 0.7586206896551724
$ws.Cells.Item(11, 14).Value = 0.7719298245614034
$ws.Cells.Item(11, 15).Value = 0.003896103896103842
$ws.Cells.Item(11, 16).Value = 0.01724137931034475
$ws.Cells.Item(11, 17).Value = 0.01086787765874853
$ws.Cells.Item(11, 18).Value = 0.004983388704318867
$ws.Cells.Item(11, 19).Value = 0.02325581395348827
$ws.Cells.Item(11, 20).Value = 0.01427988576091376
# Row 12
$ws.Cells.Item(12, 5).Value = 0.5232558139534884
$ws.Cells.Item(12, 6).Value = 0.7758620689655172
$ws.Cells.Item(12, 7).Value = 0.625
$ws.Cells.Item(12, 8).Value = 0.5232558139534884
$ws.Cells.Item(12, 9).Value = 0.7758620689655172
$ws.Cells.Item(12, 10).Value = 0.625
$ws.Cells.Item(12, 11).Value = 0.1563342318059299
$ws.Cells.Item(12, 12).Value = 0.5268817204301075
$ws.Cells.Item(12, 13).Value = 0.8448275862068966
$ws.Cells.Item(12, 14).Value = 0.6490066225165563
$ws.Cells.Item(12, 15).Value = 0.00362590647661909
$ws.Cells.Item(12, 16).Value = 0.06896551724137934
$ws.Cells.Item(12, 17).Value = 0.02400662251655628
$ws.Cells.Item(12, 18).Value = 0.006929510155316482
$ws.Cells.Item(12, 19).Value = 0.08888888888888892
$ws.Cells.Item(12, 20).Value = 0.03841059602649004
# Row 13
$ws.Cells.Item(13, 5).Value = 0.5232558139534884
$ws.Cells.Item(13, 6).Value = 0.7758620689655172
$ws.Cells.Item(13, 7).Value = 0.625
$ws.Cells.Item(13, 8).Value = 0.5232558139534884
$ws.Cells.Item(13, 9).Value = 0.7758620689655172
$ws.Cells.Item(13, 10).Value = 0.625
$ws.Cells.Item(13, 11).Value = 0.1563342318059299
$ws.Cells.Item(13, 12).Value = 0.5268817204301075
$ws.Cells.Item(13, 13).Value = 0.8448275862068966
$ws.Cells.Item(13, 14).Value = 0.6490066225165563
$ws.Cells.Item(13, 15).Value = 0.00362590647661909
$ws.Cells.Item(13, 16).Value = 0.06896551724137934
$ws.Cells.Item(13, 17).Value = 0.02400662251655628
$ws.Cells.Item(13, 18).Value = 0.006929510155316482
$ws.Cells.Item(13, 19).Value = 0.08888888888888892
$ws.Cells.Item(13, 20).Value = 0.03841059602649004
# Row 14
$ws.Cells.Item(14, 5).Value = 0.5232558139534884
$ws.Cells.Item(14, 6).Value = 0.7758620689655172
$ws.Cells.Item(14, 7).Value = 0.625
$ws.Cells.Item(14, 8).Value = 0.5232558139534884
$ws.Cells.Item(14, 9).Value = 0.7758620689655172
$ws.Cells.Item(14, 10).Value = 0.625
$ws.Cells.Item(14, 11).Value = 0.1563342318059299
$ws.Cells.Item(14, 12).Value = 0.5232558139534884
$ws.Cells.Item(14, 13).Value = 0.7758620689655172
$ws.Cells.Item(14, 14).Value = 0.625
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 0
# Row 15
$ws.Cells.Item(15, 5).Value = 0.5232558139534884
$ws.Cells.Item(15, 6).Value = 0.7758620689655172
$ws.Cells.Item(15, 7).Value = 0.625
$ws.Cells.Item(15, 8).Value = 0.5232558139534884
$ws.Cells.Item(15, 9).Value = 0.7758620689655172
$ws.Cells.Item(15, 10).Value = 0.625
$ws.Cells.Item(15, 11).Value = 0.1563342318059299
$ws.Cells.Item(15, 12).Value = 0.5319148936170213
$ws.Cells.Item(15, 13).Value = 0.8620689655172413
$ws.Cells.Item(15, 14).Value = 0.6578947368421053
$ws.Cells.Item(15, 15).Value = 0.008659079663532854
$ws.Cells.Item(15, 16).Value = 0.08620689655172409
$ws.Cells.Item(15, 17).Value = 0.03289473684210531
$ws.Cells.Item(15, 18).Value = 0.0165484633569739
$ws.Cells.Item(15, 19).Value = 0.111111111111111
$ws.Cells.Item(15, 20).Value = 0.05263157894736849
# Row 16
$ws.Cells.Item(16, 5).Value = 0.5232558139534884
$ws.Cells.Item(16, 6).Value = 0.7758620689655172
$ws.Cells.Item(16, 7).Value = 0.625
$ws.Cells.Item(16, 8).Value = 0.5232558139534884
$ws.Cells.Item(16, 9).Value = 0.7758620689655172
$ws.Cells.Item(16, 10).Value = 0.625
$ws.Cells.Item(16, 11).Value = 0.1563342318059299
$ws.Cells.Item(16, 12).Value = 0.5268817204301075
$ws.Cells.Item(16, 13).Value = 0.8448275862068966
$ws.Cells.Item(16, 14).Value = 0.6490066225165563
$ws.Cells.Item(16, 15).Value = 0.00362590647661909
$ws.Cells.Item(16, 16).Value = 0.06896551724137934
$ws.Cells.Item(16, 17).Value = 0.02400662251655628
$ws.Cells.Item(16, 18).Value = 0.006929510155316482
$ws.Cells.Item(16, 19).Value = 0.08888888888888892
$ws.Cells.Item(16, 20).Value = 0.03841059602649004
# Row 17
$ws.Cells.Item(17, 5).Value = 0.5232558139534884
$ws.Cells.Item(17, 6).Value = 0.7758620689655172
$ws.Cells.Item(17, 7).Value = 0.625
$ws.Cells.Item(17, 8).Value = 0.5232558139534884
$ws.Cells.Item(17, 9).Value = 0.7758620689655172
$ws.Cells.Item(17, 10).Value = 0.625
$ws.Cells.Item(17, 11).Value = 0.1563342318059299
$ws.Cells.Item(17, 12).Value = 0.5268817204301075
$ws.Cells.Item(17, 13).Value = 0.8448275862068966
$ws.Cells.Item(17, 14).Value = 0.6490066225165563
$ws.Cells.Item(17, 15).Value = 0.00362590647661909
$ws.Cells.Item(17, 16).Value = 0.06896551724137934
$ws.Cells.Item(17, 17).Value = 0.02400662251655628
$ws.Cells.Item(17, 18).Value = 0.006929510155316482
$ws.Cells.Item(17, 19).Value = 0.08888888888888892
$ws.Cells.Item(17, 20).Value = 0.03841059602649004
# Row 18
$ws.Cells.Item(18, 5).Value = 0.5232558139534884
$ws.Cells.Item(18, 6).Value = 0.7758620689655172
$ws.Cells.Item(18, 7).Value = 0.625
$ws.Cells.Item(18, 8).Value = 0.5232558139534884
$ws.Cells.Item(18, 9).Value = 0.7758620689655172
$ws.Cells.Item(18, 10).Value = 0.625
$ws.Cells.Item(18, 11).Value = 0.1563342318059299
$ws.Cells.Item(18, 12).Value = 0.5268817204301075
$ws.Cells.Item(18, 13).Value = 0.8448275862068966
$ws.Cells.Item(18, 14).Value = 0.6490066225165563
$ws.Cells.Item(18, 15).Value = 0.00362590647661909
$ws.Cells.Item(18, 16).Value = 0.06896551724137934
$ws.Cells.Item(18, 17).Value = 0.02400662251655628
$ws.Cells.Item(18, 18).Value = 0.006929510155316482
$ws.Cells.Item(18, 19).Value = 0.08888888888888892
$ws.Cells.Item(18, 20).Value = 0.03841059602649004
# Row 19
$ws.Cells.Item(19, 5).Value = 0.5232558139534884
$ws.Cells.Item(19, 6).Value = 0.7758620689655172
$ws.Cells.Item(19, 7).Value = 0.625
$ws.Cells.Item(19, 8).Value = 0.5232558139534884
$ws.Cells.Item(19, 9).Value = 0.7758620689655172
$ws.Cells.Item(19, 10).Value = 0.625
$ws.Cells.Item(19, 11).Value = 0.1563342318059299
$ws.Cells.Item(19, 12).Value = 0.5232558139534884
$ws.Cells.Item(19, 13).Value = 0.7758620689655172
$ws.Cells.Item(19, 14).Value = 0.625
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 0
$ws.Cells.Item(19, 18).Value = 0
$ws.Cells.Item(19, 19).Value = 0
$ws.Cells.Item(19, 20).Value = 0
# Row 20
$ws.Cells.Item(20, 5).Value = 0.5232558139534884
$ws.Cells.Item(20, 6).Value = 0.7758620689655172
$ws.Cells.Item(20, 7).Value = 0.625
$ws.Cells.Item(20, 8).Value = 0.5232558139534884
$ws.Cells.Item(20, 9).Value = 0.7758620689655172
$ws.Cells.Item(20, 10).Value = 0.625
$ws.Cells.Item(20, 11).Value = 0.1563342318059299
$ws.Cells.Item(20, 12).Value = 0.5274725274725275
$ws.Cells.Item(20, 13).Value = 0.8275862068965517
$ws.Cells.Item(20, 14).Value = 0.6442953020134228
$ws.Cells.Item(20, 15).Value = 0.004216713519039073
$ws.Cells.Item(20, 16).Value = 0.05172413793103448
$ws.Cells.Item(20, 17).Value = 0.01929530201342278
$ws.Cells.Item(20, 18).Value = 0.008058608058608004
$ws.Cells.Item(20, 19).Value = 0.06666666666666665
$ws.Cells.Item(20, 20).Value = 0.03087248322147644
# Row 21
$ws.Cells.Item(21, 5).Value = 0.5232558139534884
$ws.Cells.Item(21, 6).Value = 0.7758620689655172
$ws.Cells.Item(21, 7).Value = 0.625
$ws.Cells.Item(21, 8).Value = 0.5232558139534884
$ws.Cells.Item(21, 9).Value = 0.7758620689655172
$ws.Cells.Item(21, 10).Value = 0.625
$ws.Cells.Item(21, 11).Value = 0.1563342318059299
$ws.Cells.Item(21, 12).Value = 0.5268817204301075
$ws.Cells.Item(21, 13).Value = 0.8448275862068966
$ws.Cells.Item(21, 14).Value = 0.6490066225165563
$ws.Cells.Item(21, 15).Value = 0.00362590647661909
$ws.Cells.Item(21, 16).Value = 0.06896551724137934
$ws.Cells.Item(21, 17).Value = 0.02400662251655628
$ws.Cells.Item(21, 18).Value = 0.006929510155316482
$ws.Cells.Item(21, 19).Value = 0.08888888888888892
$ws.Cells.Item(21, 20).Value = 0.03841059602649004
# Row 22
$ws.Cells.Item(22, 5).Value = 0.8333333333333334
$ws.Cells.Item(22, 6).Value = 0.5172413793103449
$ws.Cells.Item(22, 7).Value = 0.6382978723404256
$ws.Cells.Item(22, 8).Value = 0.8333333333333334
$ws.Cells.Item(22, 9).Value = 0.5172413793103449
$ws.Cells.Item(22, 10).Value = 0.6382978723404256
$ws.Cells.Item(22, 11).Value = 0.1563342318059299
$ws.Cells.Item(22, 12).Value = 0.8536585365853658
$ws.Cells.Item(22, 13).Value = 0.603448275862069
$ws.Cells.Item(22, 14).Value = 0.7070707070707071
$ws.Cells.Item(22, 15).Value = 0.02032520325203246
$ws.Cells.Item(22, 16).Value = 0.08620689655172409
$ws.Cells.Item(22, 17).Value = 0.06877283473028151
$ws.Cells.Item(22, 18).Value = 0.02439024390243896
$ws.Cells.Item(22, 19).Value = 0.1666666666666665
$ws.Cells.Item(22, 20).Value = 0.1077441077441077
# Row 23
$ws.Cells.Item(23, 5).Value = 0.8333333333333334
$ws.Cells.Item(23, 6).Value = 0.5172413793103449
$ws.Cells.Item(23, 7).Value = 0.6382978723404256
$ws.Cells.Item(23, 8).Value = 0.8333333333333334
$ws.Cells.Item(23, 9).Value = 0.5172413793103449
$ws.Cells.Item(23, 10).Value = 0.6382978723404256
$ws.Cells.Item(23, 11).Value = 0.1563342318059299
$ws.Cells.Item(23, 12).Value = 0.8536585365853658
$ws.Cells.Item(23, 13).Value = 0.603448275862069
$ws.Cells.Item(23, 14).Value = 0.7070707070707071
$ws.Cells.Item(23, 15).Value = 0.02032520325203246
$ws.Cells.Item(23, 16).Value = 0.08620689655172409
$ws.Cells.Item(23, 17).Value = 0.06877283473028151
$ws.Cells.Item(23, 18).Value = 0.02439024390243896
$ws.Cells.Item(23, 19).Value = 0.1666666666666665
$ws.Cells.Item(23, 20).Value = 0.1077441077441077
# Row 24
$ws.Cells.Item(24, 5).Value = 0.8333333333333334
$ws.Cells.Item(24, 6).Value = 0.5172413793103449
$ws.Cells.Item(24, 7).Value = 0.6382978723404256
$ws.Cells.Item(24, 8).Value = 0.8333333333333334
$ws.Cells.Item(24, 9).Value = 0.5172413793103449
$ws.Cells.Item(24, 10).Value = 0.6382978723404256
$ws.Cells.Item(24, 11).Value = 0.1563342318059299
$ws.Cells.Item(24, 12).Value = 0.8333333333333334
$ws.Cells.Item(24, 13).Value = 0.5172413793103449
$ws.Cells.Item(24, 14).Value = 0.6382978723404256
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0
$ws.Cells.Item(24, 18).Value = 0
$ws.Cells.Item(24, 19).Value = 0
$ws.Cells.Item(24, 20).Value = 0
# Row 25
$ws.Cells.Item(25, 5).Value = 0.8333333333333334
$ws.Cells.Item(25, 6).Value = 0.5172413793103449
$ws.Cells.Item(25, 7).Value = 0.6382978723404256
$ws.Cells.Item(25, 8).Value = 0.8333333333333334
$ws.Cells.Item(25, 9).Value = 0.5172413793103449
$ws.Cells.Item(25, 10).Value = 0.6382978723404256
$ws.Cells.Item(25, 11).Value = 0.1563342318059299
$ws.Cells.Item(25, 12).Value = 0.8604651162790697
$ws.Cells.Item(25, 13).Value = 0.6379310344827587
$ws.Cells.Item(25, 14).Value = 0.7326732673267328
$ws.Cells.Item(25, 15).Value = 0.02713178294573637
$ws.Cells.Item(25, 16).Value = 0.1206896551724138
$ws.Cells.Item(25, 17).Value = 0.0943753949863072
$ws.Cells.Item(25, 18).Value = 0.03255813953488364
$ws.Cells.Item(25, 19).Value = 0.2333333333333334
$ws.Cells.Item(25, 20).Value = 0.1478547854785479
# Row 26
$ws.Cells.Item(26, 5).Value = 0.8333333333333334
$ws.Cells.Item(26, 6).Value = 0.5172413793103449
$ws.Cells.Item(26, 7).Value = 0.6382978723404256
$ws.Cells.Item(26, 8).Value = 0.8333333333333334
$ws.Cells.Item(26, 9).Value = 0.5172413793103449
$ws.Cells.Item(26, 10).Value = 0.6382978723404256
$ws.Cells.Item(26, 11).Value = 0.1563342318059299
$ws.Cells.Item(26, 12).Value = 0.8536585365853658
$ws.Cells.Item(26, 13).Value = 0.603448275862069
$ws.Cells.Item(26, 14).Value = 0.7070707070707071
$ws.Cells.Item(26, 15).Value = 0.02032520325203246
$ws.Cells.Item(26, 16).Value = 0.08620689655172409
$ws.Cells.Item(26, 17).Value = 0.06877283473028151
$ws.Cells.Item(26, 18).Value = 0.02439024390243896
$ws.Cells.Item(26, 19).Value = 0.1666666666666665
$ws.Cells.Item(26, 20).Value = 0.1077441077441077
# Row 27
$ws.Cells.Item(27, 5).Value = 0.8333333333333334
$ws.Cells.Item(27, 6).Value = 0.5172413793103449
$ws.Cells.Item(27, 7).Value = 0.6382978723404256
$ws.Cells.Item(27, 8).Value = 0.8333333333333334
$ws.Cells.Item(27, 9).Value = 0.5172413793103449
$ws.Cells.Item(27, 10).Value = 0.6382978723404256
$ws.Cells.Item(27, 11).Value = 0.1563342318059299
$ws.Cells.Item(27, 12).Value = 0.8333333333333334
$ws.Cells.Item(27, 13).Value = 0.5172413793103449
$ws.Cells.Item(27, 14).Value = 0.6382978723404256
$ws.Cells.Item(27, 15).Value = 0
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(27, 17).Value = 0
$ws.Cells.Item(27, 18).Value = 0
$ws.Cells.Item(27, 19).Value = 0
$ws.Cells.Item(27, 20).Value = 0
# Row 28
$ws.Cells.Item(28, 5).Value = 0.8333333333333334
$ws.Cells.Item(28, 6).Value = 0.5172413793103449
$ws.Cells.Item(28, 7).Value = 0.6382978723404256
$ws.Cells.Item(28, 8).Value = 0.8333333333333334
$ws.Cells.Item(28, 9).Value = 0.5172413793103449
$ws.Cells.Item(28, 10).Value = 0.6382978723404256
$ws.Cells.Item(28, 11).Value = 0.1563342318059299
$ws.Cells.Item(28, 12).Value = 0.8421052631578947
$ws.Cells.Item(28, 13).Value = 0.5517241379310345
$ws.Cells.Item(28, 14).Value = 0.6666666666666666
$ws.Cells.Item(28, 15).Value = 0.00877192982456132
$ws.Cells.Item(28, 16).Value = 0.03448275862068961
$ws.Cells.Item(28, 17).Value = 0.02836879432624106
$ws.Cells.Item(28, 18).Value = 0.01052631578947358
$ws.Cells.Item(28, 19).Value = 0.06666666666666658
$ws.Cells.Item(28, 20).Value = 0.04444444444444434
# Row 29
$ws.Cells.Item(29, 5).Value = 0.8333333333333334
$ws.Cells.Item(29, 6).Value = 0.5172413793103449
$ws.Cells.Item(29, 7).Value = 0.6382978723404256
$ws.Cells.Item(29, 8).Value = 0.8333333333333334
$ws.Cells.Item(29, 9).Value = 0.5172413793103449
$ws.Cells.Item(29, 10).Value = 0.6382978723404256
$ws.Cells.Item(29, 11).Value = 0.1563342318059299
$ws.Cells.Item(29, 12).Value = 0.8421052631578947
$ws.Cells.Item(29, 13).Value = 0.5517241379310345
$ws.Cells.Item(29, 14).Value = 0.6666666666666666
$ws.Cells.Item(29, 15).Value = 0.00877192982456132
$ws.Cells.Item(29, 16).Value = 0.03448275862068961
$ws.Cells.Item(29, 17).Value = 0.02836879432624106
$ws.Cells.Item(29, 18).Value = 0.01052631578947358
$ws.Cells.Item(29, 19).Value = 0.06666666666666658
$ws.Cells.Item(29, 20).Value = 0.04444444444444434
# Row 30
$ws.Cells.Item(30, 5).Value = 0.8333333333333334
$ws.Cells.Item(30, 6).Value = 0.5172413793103449
$ws.Cells.Item(30, 7).Value = 0.6382978723404256
$ws.Cells.Item(30, 8).Value = 0.8333333333333334
$ws.Cells.Item(30, 9).Value = 0.5172413793103449
$ws.Cells.Item(30, 10).Value = 0.6382978723404256
$ws.Cells.Item(30, 11).Value = 0.1563342318059299
$ws.Cells.Item(30, 12).Value = 0.8421052631578947
$ws.Cells.Item(30, 13).Value = 0.5517241379310345
$ws.Cells.Item(30, 14).Value = 0.6666666666666666
$ws.Cells.Item(30, 15).Value = 0.00877192982456132
$ws.Cells.Item(30, 16).Value = 0.03448275862068961
$ws.Cells.Item(30, 17).Value = 0.02836879432624106
$ws.Cells.Item(30, 18).Value = 0.01052631578947358
$ws.Cells.Item(30, 19).Value = 0.06666666666666658
$ws.Cells.Item(30, 20).Value = 0.04444444444444434
# Row 31
$ws.Cells.Item(31, 5).Value = 0.8333333333333334
$ws.Cells.Item(31, 6).Value = 0.5172413793103449
$ws.Cells.Item(31, 7).Value = 0.6382978723404256
$ws.Cells.Item(31, 8).Value = 0.8333333333333334
$ws.Cells.Item(31, 9).Value = 0.5172413793103449
$ws.Cells.Item(31, 10).Value = 0.6382978723404256
$ws.Cells.Item(31, 11).Value = 0.1563342318059299
$ws.Cells.Item(31, 12).Value = 0.8333333333333334
$ws.Cells.Item(31, 13).Value = 0.5172413793103449
$ws.Cells.Item(31, 14).Value = 0.6382978723404256
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(31, 17).Value = 0
$ws.Cells.Item(31, 18).Value = 0
$ws.Cells.Item(31, 19).Value = 0
$ws.Cells.Item(31, 20).Value = 0
# Row 32
$ws.Cells.Item(32, 5).Value = 0.2666666666666666
$ws.Cells.Item(32, 6).Value = 0.0689655172413793
$ws.Cells.Item(32, 7).Value = 0.1095890410958904
$ws.Cells.Item(32, 8).Value = 0.2666666666666666
$ws.Cells.Item(32, 9).Value = 0.0689655172413793
$ws.Cells.Item(32, 10).Value = 0.1095890410958904
$ws.Cells.Item(32, 11).Value = 0.1563342318059299
$ws.Cells.Item(32, 12).Value = 0.4814814814814815
$ws.Cells.Item(32, 13).Value = 0.896551724137931
$ws.Cells.Item(32, 14).Value = 0.6265060240963856
$ws.Cells.Item(32, 15).Value = 0.2148148148148148
$ws.Cells.Item(32, 16).Value = 0.8275862068965517
$ws.Cells.Item(32, 17).Value = 0.5169169830004952
$ws.Cells.Item(32, 18).Value = 0.8055555555555558
$ws.Cells.Item(32, 19).Value = 12
$ws.Cells.Item(32, 20).Value = 4.716867469879519
# Row 33
$ws.Cells.Item(33, 5).Value = 0.2666666666666666
$ws.Cells.Item(33, 6).Value = 0.0689655172413793
$ws.Cells.Item(33, 7).Value = 0.1095890410958904
$ws.Cells.Item(33, 8).Value = 0.2666666666666666
$ws.Cells.Item(33, 9).Value = 0.0689655172413793
$ws.Cells.Item(33, 10).Value = 0.1095890410958904
$ws.Cells.Item(33, 11).Value = 0.1563342318059299
$ws.Cells.Item(33, 12).Value = 0.4814814814814815
$ws.Cells.Item(33, 13).Value = 0.896551724137931
$ws.Cells.Item(33, 14).Value = 0.6265060240963856
$ws.Cells.Item(33, 15).Value = 0.2148148148148148
$ws.Cells.Item(33, 16).Value = 0.8275862068965517
$ws.Cells.Item(33, 17).Value = 0.5169169830004952
$ws.Cells.Item(33, 18).Value = 0.8055555555555558
$ws.Cells.Item(33, 19).Value = 12
$ws.Cells.Item(33, 20).Value = 4.716867469879519
# Row 34
$ws.Cells.Item(34, 5).Value = 0.2666666666666666
$ws.Cells.Item(34, 6).Value = 0.0689655172413793
$ws.Cells.Item(34, 7).Value = 0.1095890410958904
$ws.Cells.Item(34, 8).Value = 0.2666666666666666
$ws.Cells.Item(34, 9).Value = 0.0689655172413793
$ws.Cells.Item(34, 10).Value = 0.1095890410958904
$ws.Cells.Item(34, 11).Value = 0.1563342318059299
$ws.Cells.Item(34, 12).Value = 0.3032786885245902
$ws.Cells.Item(34, 13).Value = 0.6379310344827587
$ws.Cells.Item(34, 14).Value = 0.4111111111111111
$ws.Cells.Item(34, 15).Value = 0.03661202185792356
$ws.Cells.Item(34, 16).Value = 0.5689655172413793
$ws.Cells.Item(34, 17).Value = 0.3015220700152207
$ws.Cells.Item(34, 18).Value = 0.1372950819672134
$ws.Cells.Item(34, 19).Value = 8.250000000000002
$ws.Cells.Item(34, 20).Value = 2.751388888888889
# Row 35
$ws.Cells.Item(35, 5).Value = 0.2666666666666666
$ws.Cells.Item(35, 6).Value = 0.0689655172413793
$ws.Cells.Item(35, 7).Value = 0.1095890410958904
$ws.Cells.Item(35, 8).Value = 0.2666666666666666
$ws.Cells.Item(35, 9).Value = 0.0689655172413793
$ws.Cells.Item(35, 10).Value = 0.1095890410958904
$ws.Cells.Item(35, 11).Value = 0.1563342318059299
$ws.Cells.Item(35, 12).Value = 0.4951456310679612
$ws.Cells.Item(35, 13).Value = 0.8793103448275862
$ws.Cells.Item(35, 14).Value = 0.6335403726708075
$ws.Cells.Item(35, 15).Value = 0.2284789644012946
$ws.Cells.Item(35, 16).Value = 0.8103448275862069
$ws.Cells.Item(35, 17).Value = 0.5239513315749171
$ws.Cells.Item(35, 18).Value = 0.8567961165048548
$ws.Cells.Item(35, 19).Value = 11.75
$ws.Cells.Item(35, 20).Value = 4.781055900621119
# Row 36
$ws.Cells.Item(36, 5).Value = 0.2666666666666666
$ws.Cells.Item(36, 6).Value = 0.0689655172413793
$ws.Cells.Item(36, 7).Value = 0.1095890410958904
$ws.Cells.Item(36, 8).Value = 0.2666666666666666
$ws.Cells.Item(36, 9).Value = 0.0689655172413793
$ws.Cells.Item(36, 10).Value = 0.1095890410958904
$ws.Cells.Item(36, 11).Value = 0.1563342318059299
$ws.Cells.Item(36, 12).Value = 0.4814814814814815
$ws.Cells.Item(36, 13).Value = 0.896551724137931
$ws.Cells.Item(36, 14).Value = 0.6265060240963856
$ws.Cells.Item(36, 15).Value = 0.2148148148148148
$ws.Cells.Item(36, 16).Value = 0.8275862068965517
$ws.Cells.Item(36, 17).Value = 0.5169169830004952
$ws.Cells.Item(36, 18).Value = 0.8055555555555558
$ws.Cells.Item(36, 19).Value = 12
$ws.Cells.Item(36, 20).Value = 4.716867469879519
# Row 37
$ws.Cells.Item(37, 5).Value = 0.2666666666666666
$ws.Cells.Item(37, 6).Value = 0.0689655172413793
$ws.Cells.Item(37, 7).Value = 0.1095890410958904
$ws.Cells.Item(37, 8).Value = 0.2666666666666666
$ws.Cells.Item(37, 9).Value = 0.0689655172413793
$ws.Cells.Item(37, 10).Value = 0.1095890410958904
$ws.Cells.Item(37, 11).Value = 0.1563342318059299
$ws.Cells.Item(37, 12).Value = 0.5050505050505051
$ws.Cells.Item(37, 13).Value = 0.8620689655172413
$ws.Cells.Item(37, 14).Value = 0.6369426751592357
$ws.Cells.Item(37, 15).Value = 0.2383838383838385
$ws.Cells.Item(37, 16).Value = 0.793103448275862
$ws.Cells.Item(37, 17).Value = 0.5273536340633453
$ws.Cells.Item(37, 18).Value = 0.8939393939393945
$ws.Cells.Item(37, 19).Value = 11.5
$ws.Cells.Item(37, 20).Value = 4.812101910828027
# Row 38
$ws.Cells.Item(38, 5).Value = 0.2666666666666666
$ws.Cells.Item(38, 6).Value = 0.0689655172413793
$ws.Cells.Item(38, 7).Value = 0.1095890410958904
$ws.Cells.Item(38, 8).Value = 0.2666666666666666
$ws.Cells.Item(38, 9).Value = 0.0689655172413793
$ws.Cells.Item(38, 10).Value = 0.1095890410958904
$ws.Cells.Item(38, 11).Value = 0.1563342318059299
$ws.Cells.Item(38, 12).Value = 0.5050505050505051
$ws.Cells.Item(38, 13).Value = 0.8620689655172413
$ws.Cells.Item(38, 14).Value = 0.6369426751592357
$ws.Cells.Item(38, 15).Value = 0.2383838383838385
$ws.Cells.Item(38, 16).Value = 0.793103448275862
$ws.Cells.Item(38, 17).Value = 0.5273536340633453
$ws.Cells.Item(38, 18).Value = 0.8939393939393945
$ws.Cells.Item(38, 19).Value = 11.5
$ws.Cells.Item(38, 20).Value = 4.812101910828027
# Row 39
$ws.Cells.Item(39, 5).Value = 0.2666666666666666
$ws.Cells.Item(39, 6).Value = 0.0689655172413793
$ws.Cells.Item(39, 7).Value = 0.1095890410958904
$ws.Cells.Item(39, 8).Value = 0.2666666666666666
$ws.Cells.Item(39, 9).Value = 0.0689655172413793
$ws.Cells.Item(39, 10).Value = 0.1095890410958904
$ws.Cells.Item(39, 11).Value = 0.1563342318059299
$ws.Cells.Item(39, 12).Value = 0.2941176470588235
$ws.Cells.Item(39, 13).Value = 0.08620689655172414
$ws.Cells.Item(39, 14).Value = 0.1333333333333333
$ws.Cells.Item(39, 15).Value = 0.02745098039215693
$ws.Cells.Item(39, 16).Value = 0.01724137931034485
$ws.Cells.Item(39, 17).Value = 0.02374429223744293
$ws.Cells.Item(39, 18).Value = 0.1029411764705885
$ws.Cells.Item(39, 19).Value = 0.2500000000000003
$ws.Cells.Item(39, 20).Value = 0.2166666666666667
# Row 40
$ws.Cells.Item(40, 5).Value = 0.2666666666666666
$ws.Cells.Item(40, 6).Value = 0.0689655172413793
$ws.Cells.Item(40, 7).Value = 0.1095890410958904
$ws.Cells.Item(40, 8).Value = 0.2666666666666666
$ws.Cells.Item(40, 9).Value = 0.0689655172413793
$ws.Cells.Item(40, 10).Value = 0.1095890410958904
$ws.Cells.Item(40, 11).Value = 0.1563342318059299
$ws.Cells.Item(40, 12).Value = 0.5050505050505051
$ws.Cells.Item(40, 13).Value = 0.8620689655172413
$ws.Cells.Item(40, 14).Value = 0.6369426751592357
$ws.Cells.Item(40, 15).Value = 0.2383838383838385
$ws.Cells.Item(40, 16).Value = 0.793103448275862
$ws.Cells.Item(40, 17).Value = 0.5273536340633453
$ws.Cells.Item(40, 18).Value = 0.8939393939393945
$ws.Cells.Item(40, 19).Value = 11.5
$ws.Cells.Item(40, 20).Value = 4.812101910828027
# Row 41
$ws.Cells.Item(41, 5).Value = 0.2666666666666666
$ws.Cells.Item(41, 6).Value = 0.0689655172413793
$ws.Cells.Item(41, 7).Value = 0.1095890410958904
$ws.Cells.Item(41, 8).Value = 0.2666666666666666
$ws.Cells.Item(41, 9).Value = 0.0689655172413793
$ws.Cells.Item(41, 10).Value = 0.1095890410958904
$ws.Cells.Item(41, 11).Value = 0.1563342318059299
$ws.Cells.Item(41, 12).Value = 0.5050505050505051
$ws.Cells.Item(41, 13).Value = 0.8620689655172413
$ws.Cells.Item(41, 14).Value = 0.6369426751592357
$ws.Cells.Item(41, 15).Value = 0.2383838383838385
$ws.Cells.Item(41, 16).Value = 0.793103448275862
$ws.Cells.Item(41, 17).Value = 0.5273536340633453
$ws.Cells.Item(41, 18).Value = 0.8939393939393945
$ws.Cells.Item(41, 19).Value = 11.5
$ws.Cells.Item(41, 20).Value = 4.812101910828027
